$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.675.72"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.495.38"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "587.63"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "175.68"
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  +6.31%  "
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  +4.13%  "
$ws.Range("D12").Value = "4.94"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "25.81"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.917.21"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "67.480.01"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "2.481.39"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "11.10"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("D20").Value = "352.18"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "70.53"
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").Value = "9.22"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("D27").Value = "2.625.25"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "0.0₃0915"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").Value = "509.31"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("E31").Value = "  +3.18%  "
$ws.Range("E32").Value = "  +4.11%  "
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +8.95%  "
$ws.Range("D36").Value = "163.51"
$ws.Range("E36").Value = "  +3.27%  "
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "18.67"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").Value = "4.89"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("D45").Value = "146.09"
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.517"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0257"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").Value = "0.0746"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("D51").Value = "0.587"
$ws.Range("E51").Value = "  +1.54%  "
